$d = $word.ActiveDocument
$apos = [char]0x2019

# The last paragraph currently ends with "...account information." and
# carries a leftover "_GoBack" bookmark (an editing-position marker) right
# at its very end. The edit splits this into two paragraphs: the existing
# one is left untouched (minus the bookmark), and a new bullet paragraph
# is appended with the sentence "Taxi's code identifier is showed
# somewhere on the car." The "_GoBack" bookmark ends up sitting between
# the two halves of that new sentence, mirroring where the author's
# cursor was after typing the first chunk and pausing there.

# Step 1: drop the existing _GoBack bookmark - it will be recreated in
# its new location once the new paragraph/text exist.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Step 2: append a brand-new paragraph after the last paragraph in the
# document, inheriting the same list formatting/style (bullet list,
# Georgia font) via the normal "press Enter at end of paragraph" flow.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$tail = $lastPara.Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$startPos = $newPara.Range.Start

# Step 3: type the full new sentence into the new (currently empty)
# paragraph in one go.
$chunk1 = "Taxi${apos}s code identifier is showed"
$chunk2 = " somewhere on the car."
$newPara.Range.InsertBefore($chunk1 + $chunk2)

# Step 4: drop a fresh _GoBack bookmark exactly at the boundary between
# the two chunks, matching the author's edit history.
$bmPos = $startPos + $chunk1.Length
$bmRng = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRng)
